$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new literal text value. Updated crypto price/volume symbols
# (GitHub Actions scheduled refresh), plus a KickToken/CEJI row swap.
$updates = @(
    @{ Cell = 'D2'; Value = '276.15' }
    @{ Cell = 'D3'; Value = '21.07' }
    @{ Cell = 'D4'; Value = '6.212' }
    @{ Cell = 'D5'; Value = '0.06186' }
    @{ Cell = 'D6'; Value = '3.580' }
    @{ Cell = 'D7'; Value = '1.526' }
    @{ Cell = 'D8'; Value = '6.544' }
    @{ Cell = 'D10'; Value = '0.1642' }
    @{ Cell = 'D11'; Value = '0.08211' }
    @{ Cell = 'D12'; Value = '0.03440' }
    @{ Cell = 'D13'; Value = '0.03127' }
    @{ Cell = 'D15'; Value = '3.773' }
    @{ Cell = 'D16'; Value = '0.001616' }
    @{ Cell = 'D17'; Value = '0.04699' }
    @{ Cell = 'D18'; Value = '0.006448' }
    @{ Cell = 'D19'; Value = '0.006139' }
    @{ Cell = 'D22'; Value = '3.728' }
    @{ Cell = 'D24'; Value = '0.01386' }
    @{ Cell = 'D28'; Value = '0.0002738' }
    @{ Cell = 'D40'; Value = '0.04668' }
    @{ Cell = 'B41'; Value = 'CEJI' }
    @{ Cell = 'C41'; Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji' }
    @{ Cell = 'D41'; Value = '0.007503' }
    @{ Cell = 'E41'; Value = '40CEJICEJIBestin24h' }
    @{ Cell = 'B42'; Value = 'KickToken' }
    @{ Cell = 'C42'; Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick' }
    @{ Cell = 'D42'; Value = '0.007030' }
    @{ Cell = 'E42'; Value = '41KickTokenKICK' }
    @{ Cell = 'D43'; Value = '0.1105' }
    @{ Cell = 'D45'; Value = '0.00006297' }
    @{ Cell = 'D47'; Value = '0.8455' }
    @{ Cell = 'D49'; Value = '0.00001901' }
)

foreach ($u in $updates) {
    $c = $ws.Range($u.Cell)
    # Force text storage (many values are numeric-looking strings like
    # "276.15" which Excel would otherwise auto-convert to a Number).
    $c.NumberFormat = "@"
    $c.Value = $u.Value
    # Restore the default style so no stray formatting diff is left behind.
    $c.Style = "Normal"
}

